$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Build the new "plain body" cell style (no fill, theme-coloured Arial text)
# on an unused scratch cell first, then fan it out (format-only) to the
# cells that need it, before any values are (re)written. Doing this before
# touching row 2-4 values keeps PasteSpecial(formats) from ever creating a
# phantom D3 cell (D3 has no value in either before/after state).
# ---------------------------------------------------------------------------

$scratch = $ws.Range("J10")
$scratch.Font.ThemeColor = 1
$scratch.Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("A3:C3").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("A4:H4").PasteSpecial(-4122)
$scratch.Clear()

# ---------------------------------------------------------------------------
# Row 1 headers: add ProductName / Country / SearchCountry, plus a run of
# styled-but-empty filler cells I1:Z1 (mirrors the source sheet's layout).
# ---------------------------------------------------------------------------

# F1 / G1 already carry the "empty header" style; just fill in their text and
# restyle them to match the bold header style used by A1:E1.
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Range("F1").Value = "ProductName"
$ws.Range("G1").Value = "Country"

# H1 "SearchCountry" uses a new style: same bold/gray-fill header look, but
# with the font colour switched to the theme's text colour instead of plain
# black.
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "SearchCountry"
$ws.Range("H1").Font.ThemeColor = 1

# I1:Z1 are empty cells that still carry a (new) style: gray header fill with
# the plain (non-bold) theme body font.
$ws.Range("A1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Font.ThemeColor = 1
$ws.Range("I1").Font.Bold = $false
$ws.Range("I1").Font.Name = "Arial"
$ws.Range("I1").Copy()
$ws.Range("J1:Z1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# New row 4: TC03 / Pabla / ... using the plain body style prepared above.
# ---------------------------------------------------------------------------

$ws.Range("A4").Value = "TC03"
$ws.Range("B4").Value = "Pabla"
$ws.Range("C4").Value = "p.perez@gmail.com"
$ws.Range("D4").Value = "password"
$ws.Range("E4").Value = "Female"
$ws.Range("F4").Value = "Nokia Edge"
$ws.Range("G4").Value = "United States of America"
$ws.Range("H4").Value = "uni"
